$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.973.81"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.253.61"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "270.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.67%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "92.25"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +14.17%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.629"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "

# Row 8
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.12%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "46.08"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.84%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0973"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.36%  "

# Row 12
$ws.Range("E12").Value = "  +19.76%  "

# Row 13
$ws.Range("E13").Value = "  +1.92%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.584.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.11"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +6.38%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.253.37"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.63%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.809"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.45%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.906.17"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0000106"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "

# Row 20
$ws.Range("E20").Value = "  +3.10%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "70.93"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.48%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "234.71"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.05"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.51"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.29%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.51"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +13.05%  "

# Row 28
$ws.Range("E28").Value = "  +5.15%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "41.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.68%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "172.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0921"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.84%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "21.02"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.52%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.93%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "

# Row 36
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0351"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("E38").Value = "  -3.29%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.51"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +26.62%  "

# Row 40
$ws.Range("E40").Value = "  +15.62%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "12.92"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.19%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "63.69"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.39"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0997"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.40"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.55%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "100.45"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("E48").Value = "  +5.28%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.20"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.441"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.474.63"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
